$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hypergolic fuel tech-tier values (column C) feeding into the
# @TechRequired = bases<tier> formula in column E.
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 10

# Widen column A to fit the longer part names.
$ws.Columns("A").ColumnWidth = 13.5

# Leave the selection on the newly populated tier column.
$ws.Range("E2:E5").Select() | Out-Null
